# Update the "three-digit ÷ one-digit" answer table: each populated cell's
# division fact is replaced by a new one. Most cells keep their position in
# the table, so a straight Find/Replace on the old literal text is enough.
# The one exception is the table's last populated row, where the new set of
# five answers lands in a different cell order than the old one (one old
# cell's text is kept but shifts position, and three cells are dropped while
# effectively three new ones appear) -- for that row we set each cell's text
# directly by (row, column) so the final left-to-right order matches exactly.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

# Simple one-for-one replacements (cell position unchanged).
Replace-Text "965÷3=321, 2" "684÷5=136, 4"
Replace-Text "647÷9=71, 8" "423÷6=70, 3"
Replace-Text "964÷5=192, 4" "867÷7=123, 6"
Replace-Text "931÷8=116, 3" "402÷4=100, 2"
Replace-Text "671÷6=111, 5" "707÷5=141, 2"

Replace-Text "343÷3=114, 1" "200÷4=50, 0"
Replace-Text "380÷9=42, 2" "754÷6=125, 4"
Replace-Text "208÷8=26, 0" "653÷7=93, 2"
Replace-Text "141÷9=15, 6" "726÷9=80, 6"
Replace-Text "883÷6=147, 1" "985÷5=197, 0"

Replace-Text "457÷2=228, 1" "721÷3=240, 1"
Replace-Text "909÷7=129, 6" "978÷3=326, 0"
Replace-Text "492÷9=54, 6" "549÷4=137, 1"
Replace-Text "147÷2=73, 1" "218÷9=24, 2"
Replace-Text "446÷7=63, 5" "755÷8=94, 3"

Replace-Text "208÷6=34, 4" "856÷7=122, 2"
Replace-Text "771÷4=192, 3" "372÷3=124, 0"
Replace-Text "547÷7=78, 1" "989÷5=197, 4"
Replace-Text "400÷3=133, 1" "899÷8=112, 3"
Replace-Text "795÷2=397, 1" "915÷7=130, 5"

# Last populated row (table row 17): the five answers are replaced with a
# new set of five, and "594÷9=66, 0" moves from column 2 to column 5, so we
# set each cell explicitly rather than relying on Find/Replace.
$t = $d.Tables.Item(1)
$t.Cell(17, 1).Range.Text = "389÷6=64, 5"
$t.Cell(17, 2).Range.Text = "376÷5=75, 1"
$t.Cell(17, 3).Range.Text = "860÷3=286, 2"
$t.Cell(17, 4).Range.Text = "896÷9=99, 5"
$t.Cell(17, 5).Range.Text = "594÷9=66, 0"
